# FPVTL-1970: split the {{respondents}} merge field into its own run (with
# gramStart/gramEnd proof-error markers matching Word's grammar-check
# wrapping around "}}The"), and populate the previously-empty trailing
# paragraph with the new {{futureHearingClause}} merge field.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split "{{respondents}}The [ordinal] respondent is [name], the
#    [relationship]" into three runs:
#      "{{respondents" | "}}The" | " [ordinal] respondent is [name], the [relationship]"
#    wrapped with <w:proofErr w:type="gramStart"/> / gramEnd around the
#    "}}The" run - leaving the rest of the paragraph (the trailing
#    " [representativeClause]" runs) untouched.
# ---------------------------------------------------------------------

$respondentsParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*{{respondents}}The*") {
        $respondentsParaIndex = $i
        break
    }
}
if ($respondentsParaIndex -eq -1) {
    throw "Could not locate the paragraph containing {{respondents}}The ..."
}

$respondentsParagraph = $d.Paragraphs.Item($respondentsParaIndex)
$respondentsRange = $respondentsParagraph.Range

$respondentsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="0898C088" w14:textId="175DAA7A" w:rsidR="000C0DAE" w:rsidRDefault="000C0DAE" w:rsidP="00563DF5">
<w:r w:rsidRPr="000C0DAE"><w:t>{{respondents</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>}}The</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> [ordinal] respondent is [name], the [relationship]</w:t></w:r>
<w:r w:rsidR="00A462E6"><w:t xml:space="preserve"> </w:t></w:r>
<w:r w:rsidRPr="000C0DAE"><w:t>[</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r w:rsidR="00680203" w:rsidRPr="00680203"><w:t>representativeClause</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r w:rsidRPr="000C0DAE"><w:t>]</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$null = $respondentsRange.InsertXML($respondentsXml)

# ---------------------------------------------------------------------
# 2) Turn the last (empty) paragraph of the document into
#    "{{futureHearingClause}}" (spell-checked merge field name).
#
#    InsertXML on the document's very last paragraph inserts a sibling
#    paragraph just *before* it instead of replacing it in place (the
#    final paragraph mark is a structural sentinel), so: insert the new
#    paragraph there, then delete the now-redundant empty paragraph that
#    gets pushed after it.
# ---------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastParagraph = $d.Paragraphs.Item($lastIndex)
$lastRange = $lastParagraph.Range

$futureHearingXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="08676980" w14:textId="77777777" w:rsidR="009F19C0" w:rsidRDefault="009F19C0">
<w:r><w:t>{{</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>futureHearingClause</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>}}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$null = $lastRange.InsertXML($futureHearingXml)

$newCount = $d.Paragraphs.Count
if ($newCount -gt $lastIndex) {
    $danglingParagraph = $d.Paragraphs.Item($newCount)
    $danglingRange = $d.Range($danglingParagraph.Range.Start - 1, $danglingParagraph.Range.End)
    $null = $danglingRange.Delete()
}
